$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Add the new "2022" column (column S) --------------------------------
# Row 2 (thin header spacer row): S2 should look like its neighbour R2
$ws.Range("R2").Copy() | Out-Null
$ws.Range("S2").PasteSpecial(-4122) | Out-Null

# Row 3 (year headers): S3 should look like its neighbour R3, with value 2022
$ws.Range("R3").Copy() | Out-Null
$ws.Range("S3").PasteSpecial(-4122) | Out-Null
$ws.Range("S3").Value = 2022

# Row 4 (GVA share %%): S4 should look like its neighbour R4
$ws.Range("R4").Copy() | Out-Null
$ws.Range("S4").PasteSpecial(-4122) | Out-Null
$ws.Range("S4").Value = 13.6

# Row 5 (GVA per capita): S5 should look like its neighbour R5
$ws.Range("R5").Copy() | Out-Null
$ws.Range("S5").PasteSpecial(-4122) | Out-Null
$ws.Range("S5").Value = 20

$ws.Range("A1").Select() | Out-Null
$excel.CutCopyMode = 0

# --- Update existing 2019-2021 figures in rows 4 and 5 --------------------
$ws.Range("P4").Value = 13.7
$ws.Range("Q4").Value = 13.1
$ws.Range("R4").Value = 11.8

$ws.Range("P5").Value = 13.6
$ws.Range("Q5").Value = 12.5
$ws.Range("R5").Value = 13.5

# --- Update selection shown in the saved view ------------------------------
$ws.Range("S2").Select() | Out-Null
